# Convert the position values in row 3 (columns A-M) of every test sheet
# from inches to meters (1 in = 0.0254 m).

$wb = $excel.ActiveWorkbook

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M")

foreach ($ws in $wb.Worksheets) {
    foreach ($col in $cols) {
        $cell = $ws.Range($col + "3")
        $inches = $cell.Value()
        $cell.Value = $inches * 0.0254
    }
}

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# Re-create the author's navigation/selection trail while reviewing the
# converted row on each sheet, finishing back on sheet 1.
$null = $ws2.Activate()
$null = $ws2.Range("A3:M3").EntireRow.Select()

$null = $ws3.Activate()
$null = $ws3.Range("A3:M3").EntireRow.Select()

$null = $ws4.Activate()
$null = $ws4.Range("A3:M3").EntireRow.Select()

$null = $ws1.Activate()
$null = $ws1.Range("D12").Select()
